$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected; unprotect to allow value updates, re-protect at the end
$ws.Unprotect("D382")

# Update the confidential disclosure date from 2021-05-10 to 2021-05-11
$ws.Range("A41").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-11 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for holdings rows 2-38
$ws.Range("D2").Value = 0.03183987463494551
$ws.Range("E2").Value = -0.01144381345723622
$ws.Range("D3").Value = 0.02814256048104613
$ws.Range("E3").Value = 0.01211982143602164
$ws.Range("D4").Value = 0.02825351512936099
$ws.Range("E4").Value = 0.01105408606395586
$ws.Range("D5").Value = 0.06243380880107143
$ws.Range("E5").Value = 0.01047488003410146
$ws.Range("D6").Value = 0.01607374746487219
$ws.Range("E6").Value = -0.01964937910883868
$ws.Range("D7").Value = 0.01513914534784967
$ws.Range("E7").Value = -0.004420660772452289
$ws.Range("D8").Value = 0.02966070185957643
$ws.Range("E8").Value = -0.001266724724883095
$ws.Range("D9").Value = 0.03467968743273253
$ws.Range("E9").Value = -0.008802618214648561
$ws.Range("D10").Value = 0.02857170253352318
$ws.Range("E10").Value = 0.01723889950481827
$ws.Range("D11").Value = 0.03007144931427288
$ws.Range("E11").Value = 0.008427094246799882
$ws.Range("D12").Value = 0.0111414513259374
$ws.Range("E12").Value = 0.0194959163958901
$ws.Range("D13").Value = 0.01466245130323769
$ws.Range("E13").Value = 0.002989536621823552
$ws.Range("D14").Value = 0.01423233081470848
$ws.Range("E14").Value = 0.02127045235803648
$ws.Range("D15").Value = 0.009294359746671166
$ws.Range("E15").Value = -0.03663466397170279
$ws.Range("D16").Value = 0.008104972881666379
$ws.Range("E16").Value = -0.01760104302477183
$ws.Range("D17").Value = 0.02929887620748794
$ws.Range("E17").Value = 0.007807751648043482
$ws.Range("D18").Value = 0.02567788006565698
$ws.Range("E18").Value = -0.0004191466174869873
$ws.Range("D19").Value = 0.03264943262450208
$ws.Range("E19").Value = -0.0005394228175852334
$ws.Range("D20").Value = 0.02993720788791664
$ws.Range("E20").Value = 0.001830244795241187
$ws.Range("D21").Value = 0.04582329131798468
$ws.Range("E21").Value = -0.0140498620636641
$ws.Range("D22").Value = 0.03607298037100731
$ws.Range("E22").Value = -0.02376044266030164
$ws.Range("D23").Value = 0.0333764106112564
$ws.Range("E23").Value = -0.03066369606003738
$ws.Range("D24").Value = 0.03178410377997243
$ws.Range("E24").Value = -0.02233058126004317
$ws.Range("D25").Value = 0.01478162481439069
$ws.Range("E25").Value = -0.0315739306564996
$ws.Range("D26").Value = 0.01476225178055793
$ws.Range("E26").Value = -0.005806092420265641
$ws.Range("D27").Value = 0.03146395950370592
$ws.Range("E27").Value = -0.01334063077238823
$ws.Range("D28").Value = 0.03164340467567194
$ws.Range("E28").Value = -0.02072923366150492
$ws.Range("D29").Value = 0.02902197880472687
$ws.Range("E29").Value = -0.003843353022089158
$ws.Range("D30").Value = 0.02936090905319483
$ws.Range("E30").Value = 0.005278592375366653
$ws.Range("D31").Value = 0.03349949786661804
$ws.Range("E31").Value = 0.00283896745702128
$ws.Range("D32").Value = 0.0315823502660172
$ws.Range("E32").Value = -0.007292802617230087
$ws.Range("D33").Value = 0.02857091978468145
$ws.Range("E33").Value = 0.004808120381088132
$ws.Range("D34").Value = 0.03294511599946461
$ws.Range("E34").Value = -0.02183468364655849
$ws.Range("D35").Value = 0.03039100652890809
$ws.Range("E35").Value = -0.0001159017153454789
$ws.Range("D36").Value = 0.03154869206582292
$ws.Range("E36").Value = -0.01662324773601298
$ws.Range("D37").Value = 0.03350634691898315
$ws.Range("E37").Value = -0.02238004018128292
$ws.Range("E38").Value = -0.005514661277179811

# Restore sheet protection
$ws.Protect("D382")
